$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "challenges": row 3 (G2) is dropped; row 2 (G1) is updated with
# some of row 3's values (K,L) while B/D get "0"/"G0" and M2 is cleared.
# ----------------------------------------------------------------------
$wsChallenges = $wb.Worksheets.Item("challenges")

$wsChallenges.Range("B2").Value = 0
$wsChallenges.Range("D2").Value = "G0"
$wsChallenges.Range("K2").Value = 0
$wsChallenges.Range("L2").Value = 25
$wsChallenges.Range("M2").ClearContents()

# Delete row 3 entirely (shifts nothing up below it, it's the last row)
$wsChallenges.Rows.Item(3).Delete()

# ----------------------------------------------------------------------
# Sheet "tasks": rows 2-4 are each overwritten with the content that used
# to live one row below them, and the old last row (5) is removed.
# ----------------------------------------------------------------------
$wsTasks = $wb.Worksheets.Item("tasks")

# Row 2 becomes old row 3 ("Buy half heart"), minus F/G, A becomes 0
$wsTasks.Range("A2").Value = 0
$wsTasks.Range("B2").Value = "Buy half heart"
$wsTasks.Range("F2").ClearContents()
$wsTasks.Range("G2").ClearContents()
$wsTasks.Range("I2").Value = "ConfusingArrowsData"
$wsTasks.Range("J2").Value = "ConfusingArrowsData"
$wsTasks.Range("L2").Value = "[MINIGAME_BUY_HALF_HEART, STRICTLY_GREATER, 0],[MINIGAMESTATE_ID, EQUAL, 1]"
$wsTasks.Range("M2").Value = -5

# Row 3 becomes old row 4 ("Score 5 points"), A becomes 0
$wsTasks.Range("A3").Value = 0
$wsTasks.Range("B3").Value = "Score 5 points"
$wsTasks.Range("L3").Value = "[MINIGAME_SCORE, STRICTLY_GREATER, 5],[MINIGAMESTATE_ID, EQUAL, 1]"
$wsTasks.Range("M3").Value = 10

# Row 4 becomes old row 5 ("Walk 500 meters"), A becomes 0
$wsTasks.Range("A4").Value = 0
$wsTasks.Range("B4").Value = "Walk 500 meters"
$wsTasks.Range("I4").Value = "WALK"
$wsTasks.Range("J4").Value = "WALK"
$wsTasks.Range("L4").Value = "[DISTANCE, STRICTLY_GREATER, 499]"
$wsTasks.Range("M4").Value = 20

# Delete old row 5 (now duplicated content, last row)
$wsTasks.Rows.Item(5).Delete()
